$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / account holder info
$ws.Range("C2").Value = "Hartmut"
# Card number is a 16-digit value that must stay text (storing as a real
# Number would both lose precision and round-trip through scientific
# notation), so force text via the quote-prefix the same way typing it
# into Excel would.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 07.03.2025"

# Row 6
$ws.Range("B6").Value = "08.03."
$ws.Range("C6").Value = "09.03."
$ws.Range("D6").Value = "ZALANDO MKTPLC EU NYHWMA"
$ws.Range("E6").Value = "44,67-"

# Row 7
$ws.Range("B7").Value = "11.03."
$ws.Range("C7").Value = "12.03."
$ws.Range("D7").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E7").Value = "77,68-"

# Row 8
$ws.Range("B8").Value = "13.03."
$ws.Range("C8").Value = "14.03."
$ws.Range("D8").Value = "BURGER KING Wittmund"
$ws.Range("E8").Value = "13,80-"

# Row 9 (previously blank, now filled)
$ws.Range("B9").Value = "15.03."
$ws.Range("C9").Value = "16.03."
$ws.Range("D9").Value = "KARTENZAHLUNG JET TANKSTELLE"
$ws.Range("E9").Value = "88,57-"
$ws.Range("E9").HorizontalAlignment = -4152
$ws.Range("E9").VerticalAlignment = -4107
$ws.Range("E9").WrapText = $false

# Row 10 (previously blank, now filled)
$ws.Range("B10").Value = "16.03."
$ws.Range("C10").Value = "17.03."
$ws.Range("D10").Value = "RECHNUNG VODAFONE GMBH 50467717"
$ws.Range("E10").Value = "40,52-"
$ws.Range("E10").HorizontalAlignment = -4152
$ws.Range("E10").VerticalAlignment = -4107
$ws.Range("E10").WrapText = $false

# Row 11 (previously blank, now filled)
$ws.Range("B11").Value = "20.03."
$ws.Range("C11").Value = "21.03."
$ws.Range("D11").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 18385497"
$ws.Range("E11").Value = "83,01-"
$ws.Range("E11").HorizontalAlignment = -4152
$ws.Range("E11").VerticalAlignment = -4107
$ws.Range("E11").WrapText = $false

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 22.03.2025"
$ws.Range("E12").Value = "348,25-"

# Next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 29.03.2025"
